# Modified Work Order templates and other fixes.

$wb = $excel.ActiveWorkbook

# "Create WO" is the first / active sheet (sheet1.xml)
$ws = $wb.Worksheets.Item("Create WO")

# Update the shared string text used by B2 ("Pro-Lot Track (Lot Track)" -> "Pro-SYDATA1 (Lot track)")
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Update the saved cell selection/active cell for the sheet (was C3 / A3:XFD3 -> now B3 / B3)
$ws.Range("B3").Select()

# Add page setup (portrait orientation) to the sheet, matching the other sheet's print settings
$ws.PageSetup.Orientation = 1
